$wb = $excel.ActiveWorkbook

# --- SES scales sheet: fix curly quotes -> straight quotes in two description cells ---
$wsSes = $wb.Worksheets.Item("SES scales ")
[void]$wsSes.Activate()
$wsSes.Range("D1").Value = "'Better off' or 'poor' based on basic needs being met (e.g. water, electricity)."
$wsSes.Range("D4").Value = "1=10th grade or less, 2=vocational education, 3=high school education, 4=short tertiary education, 5=bachelor's degree or equivalent, 6=master's degree or higher."
[void]$wsSes.Range("D17").Select()

# --- Metadata sheet: clear the trailing unused/blank styled columns R:T ---
$wsMeta = $wb.Worksheets.Item("Metadata")
[void]$wsMeta.Activate()
$wsMeta.Range("R1:T1").Clear()

# --- Dictionary sheet: leave selection on A1:A17, not the active tab ---
$wsDict = $wb.Worksheets.Item("Dictionary")
[void]$wsDict.Activate()
[void]$wsDict.Range("A1:A17").Select()

# --- Final state: Metadata is the active sheet with G24 selected ---
[void]$wsMeta.Activate()
[void]$wsMeta.Range("G24").Select()
